$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Ngf", "Sorcs3", "FAPs", 1, 0.3333333333333333, 0.2000683333333333, 0.600205, 0.01182207550347303, 0.01182207550347303, 1, 0.3333333333333333, 0.004862, 0.014586, 0.02348134339170667, 0.02348134339170667, 0.0009727322366666666, 0.00875459013, 0.0002775982144997337, 0.0002775982144997337),
    @("ECs", "Ngf", "Sorcs3", "sCs", 1, 0.3333333333333333, 0.2000683333333333, 0.600205, 0.01182207550347303, 0.01182207550347303, 3, 1, 0.202196, 0.606588, 0.9765186566082933, 0.9765186566082934, 0.04045301672666667, 0.36407715054, 0.01154447728897329, 0.01154447728897329),
    @("FAPs", "Ngf", "Sorcs3", "FAPs", 1, 0.3333333333333333, 0.44525, 1.33575, 0.02630990637159653, 0.02630990637159653, 1, 0.3333333333333333, 0.004862, 0.014586, 0.02348134339170667, 0.02348134339170667, 0.0021648055, 0.0194832495, 0.0006177919461151095, 0.0006177919461151095),
    @("FAPs", "Ngf", "Sorcs3", "sCs", 1, 0.3333333333333333, 0.44525, 1.33575, 0.02630990637159653, 0.02630990637159653, 3, 1, 0.202196, 0.606588, 0.9765186566082933, 0.9765186566082934, 0.09002776900000001, 0.810249921, 0.02569211442548142, 0.02569211442548142),
    @("sCs", "Ngf", "Sorcs3", "FAPs", 3, 1, 16.277965, 48.833895, 0.9618680181249305, 0.9618680181249305, 1, 0.3333333333333333, 0.004862, 0.014586, 0.02348134339170667, 0.02348134339170667, 0.07914346583, 0.71229119247, 0.02258595323109183, 0.02258595323109183),
    @("sCs", "Ngf", "Sorcs3", "sCs", 3, 1, 16.277965, 48.833895, 0.9618680181249305, 0.9618680181249305, 3, 1, 0.202196, 0.606588, 0.9765186566082933, 0.9765186566082934, 3.29133941114, 29.62205470026, 0.9392820648938386, 0.9392820648938387)
)

$startRow = 2
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value = $rowVals[$c]
    }
}
